$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 403997.8
$ws.Range("J46").Value = 666996.3
$ws.Range("L46").Value = 2000988.9
$ws.Range("N46").Value = -2001226.9
$ws.Range("H60").Value = 403997.8
$ws.Range("J60").Value = 666996.3
$ws.Range("L60").Value = 2000988.9
$ws.Range("N60").Value = -2001956.9
$ws.Range("H62").Value = 4920
$ws.Range("I62").Value = 4800
$ws.Range("K62").Value = 4800
$ws.Range("M62").Value = -4176
$ws.Range("H65").Value = 4920
$ws.Range("I65").Value = 4800
$ws.Range("K65").Value = 24000
$ws.Range("M65").Value = -20880
$ws.Range("H74").Value = 13468
$ws.Range("I74").Value = 12702.167
$ws.Range("K74").Value = 12702.167
$ws.Range("M74").Value = -11766.167
$ws.Range("H77").Value = 13468
$ws.Range("I77").Value = 12702.167
$ws.Range("K77").Value = 63510.835
$ws.Range("M77").Value = -58830.835
$ws.Range("H106").Value = 6399.8
$ws.Range("I106").Value = 3999.5
$ws.Range("K106").Value = 3999.5
$ws.Range("M106").Value = -3368.5
$ws.Range("H111").Value = 47843
$ws.Range("I111").Value = 53611.6
$ws.Range("K111").Value = 160834.8
$ws.Range("M111").Value = -157767.8
$ws.Range("H118").Value = 1088.375
$ws.Range("I118").Value = 1101
$ws.Range("J118").Value = 1000
$ws.Range("K118").Value = 3303
$ws.Range("L118").Value = 3000
$ws.Range("M118").Value = -1646
$ws.Range("N118").Value = -6314
$ws.Range("H137").Value = 8278.68
$ws.Range("I137").Value = 3348.25
$ws.Range("J137").Value = 28000.4
$ws.Range("K137").Value = 10044.75
$ws.Range("L137").Value = 84001.20000000001
$ws.Range("M137").Value = -7494.75
$ws.Range("N137").Value = -89101.20000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7815.607
$ws.Range("I32").Value = 6214.7075
$ws.Range("K32").Value = 6214.7075
$ws.Range("M32").Value = -5927.7075
$ws.Range("H45").Value = 5451.2856
$ws.Range("I45").Value = 4776.5
$ws.Range("K45").Value = 4776.5
$ws.Range("M45").Value = -4399.5
$ws.Range("H61").Value = 5010
$ws.Range("I61").Value = 5244.827
$ws.Range("J61").Value = 3653.2222
$ws.Range("K61").Value = 5244.827
$ws.Range("L61").Value = 3653.2222
$ws.Range("M61").Value = -5032.827
$ws.Range("N61").Value = -4077.2222
$ws.Range("H74").Value = 3261.8333
$ws.Range("J74").Value = 5989.7
$ws.Range("L74").Value = 5989.7
$ws.Range("N74").Value = -7737.7
$ws.Range("H77").Value = 3261.8333
$ws.Range("J77").Value = 5989.7
$ws.Range("L77").Value = 29948.5
$ws.Range("N77").Value = -38684.5
$ws.Range("H124").Value = 28750
$ws.Range("J124").Value = 28750
$ws.Range("L124").Value = 28750
$ws.Range("N124").Value = -38570
$ws.Range("H136").Value = 5010
$ws.Range("I136").Value = 5244.827
$ws.Range("J136").Value = 3653.2222
$ws.Range("K136").Value = 15734.481
$ws.Range("L136").Value = 10959.6666
$ws.Range("M136").Value = -13184.481
$ws.Range("N136").Value = -16059.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3709.5789
$ws.Range("I20").Value = 3677.4285
$ws.Range("J20").Value = 3799.6
$ws.Range("K20").Value = 3677.4285
$ws.Range("L20").Value = 3799.6
$ws.Range("M20").Value = -3430.4285
$ws.Range("N20").Value = -4293.6
$ws.Range("H81").Value = 37055.6
$ws.Range("J81").Value = 37055.6
$ws.Range("L81").Value = 37055.6
$ws.Range("N81").Value = -39177.6
$ws.Range("H84").Value = 37055.6
$ws.Range("J84").Value = 37055.6
$ws.Range("L84").Value = 111166.8
$ws.Range("N84").Value = -121774.8
$ws.Range("H107").Value = 2274.8333
$ws.Range("I107").Value = 1474.875
$ws.Range("J107").Value = 3874.75
$ws.Range("K107").Value = 1474.875
$ws.Range("L107").Value = 3874.75
$ws.Range("M107").Value = 445.125
$ws.Range("N107").Value = -7714.75
$ws.Range("H134").Value = 2783.5
$ws.Range("I134").Value = 2529.1956
$ws.Range("J134").Value = 5708
$ws.Range("K134").Value = 7587.5868
$ws.Range("L134").Value = 17124
$ws.Range("M134").Value = -5052.5868
$ws.Range("N134").Value = -22194

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2343.9473
$ws.Range("I16").Value = 1931.4166
$ws.Range("K16").Value = 1931.4166
$ws.Range("M16").Value = -1644.4166
$ws.Range("H62").Value = 3765.6924
$ws.Range("J62").Value = 3993
$ws.Range("L62").Value = 3993
$ws.Range("N62").Value = -5241
$ws.Range("H65").Value = 3765.6924
$ws.Range("J65").Value = 3993
$ws.Range("L65").Value = 19965
$ws.Range("N65").Value = -26205
$ws.Range("H86").Value = 5234.12
$ws.Range("I86").Value = 4234.625
$ws.Range("J86").Value = 7011
$ws.Range("K86").Value = 4234.625
$ws.Range("L86").Value = 7011
$ws.Range("M86").Value = -3111.625
$ws.Range("N86").Value = -9257
$ws.Range("H89").Value = 5234.12
$ws.Range("I89").Value = 4234.625
$ws.Range("J89").Value = 7011
$ws.Range("K89").Value = 21173.125
$ws.Range("L89").Value = 35055
$ws.Range("M89").Value = -15557.125
$ws.Range("N89").Value = -46287
$ws.Range("H102").Value = 75088.664
$ws.Range("J102").Value = 75088.664
$ws.Range("L102").Value = 75088.664
$ws.Range("N102").Value = -79956.664
$ws.Range("H107").Value = 488.05713
$ws.Range("I107").Value = 316.58334
$ws.Range("K107").Value = 316.58334
$ws.Range("M107").Value = 1603.41666
$ws.Range("H113").Value = 2343.9473
$ws.Range("I113").Value = 1931.4166
$ws.Range("K113").Value = 1931.4166
$ws.Range("M113").Value = 238.5834
$ws.Range("H132").Value = 1334749.2
$ws.Range("I132").Value = 1334749.2
$ws.Range("K132").Value = 4004247.6
$ws.Range("M132").Value = -4001717.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 223684.78
$ws.Range("I131").Value = 715201.2
$ws.Range("J131").Value = 1709.6129
$ws.Range("K131").Value = 2145603.6
$ws.Range("L131").Value = 5128.8387
$ws.Range("M131").Value = -2140563.6
$ws.Range("N131").Value = -15208.8387

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 681.6316
$ws.Range("J107").Value = 1275.6666
$ws.Range("L107").Value = 1275.6666
$ws.Range("N107").Value = -5115.6666
$ws.Range("H122").Value = 1613
$ws.Range("I122").Value = 970.0909
$ws.Range("J122").Value = 5149
$ws.Range("K122").Value = 2910.2727
$ws.Range("L122").Value = 15447
$ws.Range("M122").Value = -460.2727
$ws.Range("N122").Value = -20347
$ws.Range("H132").Value = 2462.4
$ws.Range("I132").Value = 2500.4358
$ws.Range("J132").Value = 979
$ws.Range("K132").Value = 7501.307400000001
$ws.Range("L132").Value = 2937
$ws.Range("M132").Value = -4971.307400000001
$ws.Range("N132").Value = -7997

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2765.75
$ws.Range("I22").Value = 1437.8334
$ws.Range("J22").Value = 3334.8572
$ws.Range("K22").Value = 1437.8334
$ws.Range("L22").Value = 3334.8572
$ws.Range("M22").Value = -1142.8334
$ws.Range("N22").Value = -3924.8572
$ws.Range("H27").Value = 2765.75
$ws.Range("I27").Value = 1437.8334
$ws.Range("J27").Value = 3334.8572
$ws.Range("K27").Value = 1437.8334
$ws.Range("L27").Value = 3334.8572
$ws.Range("M27").Value = -1330.8334
$ws.Range("N27").Value = -3548.8572
$ws.Range("H55").Value = 347.72726
$ws.Range("I55").Value = 371.5
$ws.Range("J55").Value = 306.125
$ws.Range("K55").Value = 371.5
$ws.Range("L55").Value = 306.125
$ws.Range("M55").Value = -198.5
$ws.Range("N55").Value = -652.125
$ws.Range("H122").Value = 5251.5312
$ws.Range("I122").Value = 4657.087
$ws.Range("K122").Value = 13971.261
$ws.Range("M122").Value = -11521.261
$ws.Range("H132").Value = 3382.9768
$ws.Range("I132").Value = 3357.16
$ws.Range("J132").Value = 3418.8333
$ws.Range("K132").Value = 10071.48
$ws.Range("L132").Value = 10256.4999
$ws.Range("M132").Value = -7541.48
$ws.Range("N132").Value = -15316.4999
$ws.Range("H136").Value = 4716.689
$ws.Range("I136").Value = 4207.0347
$ws.Range("K136").Value = 12621.1041
$ws.Range("M136").Value = -10071.1041

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1325.9166
$ws.Range("I107").Value = 1862.8
$ws.Range("K107").Value = 5588.4
$ws.Range("M107").Value = -3668.4
$ws.Range("H122").Value = 2912.7144
$ws.Range("I122").Value = 2376.9
$ws.Range("K122").Value = 7130.700000000001
$ws.Range("M122").Value = -4680.700000000001
$ws.Range("H136").Value = 36332.668
$ws.Range("I136").Value = 51999.75
$ws.Range("K136").Value = 155999.25
$ws.Range("M136").Value = -153449.25

Write-Output "applied 224 cell updates across 8 sheets"